# Auto-generated Excel COM-interop edit script
# Applies the diff: updates player/team rows in "snapshot",
# and the corresponding "returned" / "new_injured" event rows,
# reflecting the 2025-12-22 scrape run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("snapshot")

$c = $ws.Range("K2")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:40:44.018337+00:00"

$c = $ws.Range("K3")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:40:46.179374+00:00"

$c = $ws.Range("K4")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:40:46.179402+00:00"

$c = $ws.Range("K5")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:40:46.179419+00:00"

$c = $ws.Range("K6")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:40:48.273600+00:00"

$c = $ws.Range("K7")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:40:49.988682+00:00"

$rng = $ws.Range("A8:K8")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "БАР"
$arr[0,1] = "Барыс"
$arr[0,2] = "barys"
$arr[0,3] = "Шил Адам"
$arr[0,4] = "41"
$arr[0,5] = "вратарь"
$arr[0,6] = "45713"
$arr[0,7] = "1369_БАР_шиладам"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/barys/team/"
$arr[0,10] = "2025-12-22T04:40:52.227021+00:00"
$rng.Value = $arr

$rng = $ws.Range("A9:K9")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "ДИН"
$arr[0,1] = "Динамо М"
$arr[0,2] = "dynamo_msk"
$arr[0,3] = "Готовец Кирилл"
$arr[0,4] = "41"
$arr[0,5] = "защитник"
$arr[0,6] = "16034"
$arr[0,7] = "1369_ДИН_готовецкирилл"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/dynamo_msk/team/"
$arr[0,10] = "2025-12-22T04:40:54.364582+00:00"
$rng.Value = $arr

$rng = $ws.Range("A10:K10")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "ЛАД"
$arr[0,1] = "Лада"
$arr[0,2] = "lada"
$arr[0,3] = "Ожгихин Алексей"
$arr[0,4] = "43"
$arr[0,5] = "нападающий"
$arr[0,6] = "23021"
$arr[0,7] = "1369_ЛАД_ожгихиналексей"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/lada/team/"
$arr[0,10] = "2025-12-22T04:40:57.756284+00:00"
$rng.Value = $arr

$rng = $ws.Range("A11:K11")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "ЛОК"
$arr[0,1] = "Локомотив"
$arr[0,2] = "lokomotiv"
$arr[0,3] = "Сергеев Андрей"
$arr[0,4] = "99"
$arr[0,5] = "защитник"
$arr[0,6] = "15416"
$arr[0,7] = "1369_ЛОК_сергеевандрей"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/lokomotiv/team/"
$arr[0,10] = "2025-12-22T04:40:59.404312+00:00"
$rng.Value = $arr

$rng = $ws.Range("A12:K12")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "ММГ"
$arr[0,1] = "Металлург Мг"
$arr[0,2] = "metallurg_mg"
$arr[0,3] = "Сиряцкий Александр"
$arr[0,4] = "74"
$arr[0,5] = "защитник"
$arr[0,6] = "42458"
$arr[0,7] = "1369_ММГ_сиряцкийалександр"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/metallurg_mg/team/"
$arr[0,10] = "2025-12-22T04:41:01.110875+00:00"
$rng.Value = $arr

$rng = $ws.Range("A13:K13")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СЕВ"
$arr[0,1] = "Северсталь"
$arr[0,2] = "severstal"
$arr[0,3] = "Ващенко Григорий"
$arr[0,4] = "16"
$arr[0,5] = "защитник"
$arr[0,6] = "14155"
$arr[0,7] = "1369_СЕВ_ващенкогригорий"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/severstal/team/"
$arr[0,10] = "2025-12-22T04:41:05.132566+00:00"
$rng.Value = $arr

$rng = $ws.Range("A14:K14")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СЕВ"
$arr[0,1] = "Северсталь"
$arr[0,2] = "severstal"
$arr[0,3] = "Смирнов Егор Д"
$arr[0,4] = "47"
$arr[0,5] = "нападающий"
$arr[0,6] = "40906"
$arr[0,7] = "1369_СЕВ_смирновегорд"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/severstal/team/"
$arr[0,10] = "2025-12-22T04:41:05.132590+00:00"
$rng.Value = $arr

$rng = $ws.Range("A15:K15")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СЕВ"
$arr[0,1] = "Северсталь"
$arr[0,2] = "severstal"
$arr[0,3] = "Фомин Макар"
$arr[0,4] = "77"
$arr[0,5] = "защитник"
$arr[0,6] = "42087"
$arr[0,7] = "1369_СЕВ_фоминмакар"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/severstal/team/"
$arr[0,10] = "2025-12-22T04:41:05.132598+00:00"
$rng.Value = $arr

$rng = $ws.Range("A16:K16")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СЕВ"
$arr[0,1] = "Северсталь"
$arr[0,2] = "severstal"
$arr[0,3] = "Шостак Константин"
$arr[0,4] = "1"
$arr[0,5] = "вратарь"
$arr[0,6] = "27876"
$arr[0,7] = "1369_СЕВ_шостакконстантин"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/severstal/team/"
$arr[0,10] = "2025-12-22T04:41:05.132606+00:00"
$rng.Value = $arr

$rng = $ws.Range("A17:K17")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СИБ"
$arr[0,1] = "Сибирь"
$arr[0,2] = "sibir"
$arr[0,3] = "Люзенков Илья"
$arr[0,4] = "86"
$arr[0,5] = "нападающий"
$arr[0,6] = "41344"
$arr[0,7] = "1369_СИБ_люзенковилья"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/sibir/team/"
$arr[0,10] = "2025-12-22T04:41:06.783156+00:00"
$rng.Value = $arr

$rng = $ws.Range("A18:K18")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СИБ"
$arr[0,1] = "Сибирь"
$arr[0,2] = "sibir"
$arr[0,3] = "Пьянов Валентин"
$arr[0,4] = "45"
$arr[0,5] = "нападающий"
$arr[0,6] = "16195"
$arr[0,7] = "1369_СИБ_пьяноввалентин"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/sibir/team/"
$arr[0,10] = "2025-12-22T04:41:06.783183+00:00"
$rng.Value = $arr

$rng = $ws.Range("A19:K19")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СИБ"
$arr[0,1] = "Сибирь"
$arr[0,2] = "sibir"
$arr[0,3] = "Чехович Иван"
$arr[0,4] = "82"
$arr[0,5] = "нападающий"
$arr[0,6] = "24581"
$arr[0,7] = "1369_СИБ_чеховичиван"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/sibir/team/"
$arr[0,10] = "2025-12-22T04:41:06.783199+00:00"
$rng.Value = $arr

$rng = $ws.Range("A20:K20")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СКА"
$arr[0,1] = "СКА"
$arr[0,2] = "ska"
$arr[0,3] = "Зайцев Никита И"
$arr[0,4] = "22"
$arr[0,5] = "защитник"
$arr[0,6] = "16024"
$arr[0,7] = "1369_СКА_зайцевникитаи"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/ska/team/"
$arr[0,10] = "2025-12-22T04:41:08.460337+00:00"
$rng.Value = $arr

$rng = $ws.Range("A21:K21")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СКА"
$arr[0,1] = "СКА"
$arr[0,2] = "ska"
$arr[0,3] = "Короткий Матвей"
$arr[0,4] = "71"
$arr[0,5] = "нападающий"
$arr[0,6] = "41566"
$arr[0,7] = "1369_СКА_короткийматвей"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/ska/team/"
$arr[0,10] = "2025-12-22T04:41:08.460358+00:00"
$rng.Value = $arr

$rng = $ws.Range("A22:K22")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СКА"
$arr[0,1] = "СКА"
$arr[0,2] = "ska"
$arr[0,3] = "Мёрфи Тревор"
$arr[0,4] = "8"
$arr[0,5] = "защитник"
$arr[0,6] = "34733"
$arr[0,7] = "1369_СКА_мерфитревор"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/ska/team/"
$arr[0,10] = "2025-12-22T04:41:08.460367+00:00"
$rng.Value = $arr

$rng = $ws.Range("A23:K23")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СОЧ"
$arr[0,1] = "ХК Сочи"
$arr[0,2] = "hc_sochi"
$arr[0,3] = "Самсонов Илья"
$arr[0,4] = "35"
$arr[0,5] = "вратарь"
$arr[0,6] = "21010"
$arr[0,7] = "1369_СОЧ_самсоновилья"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/hc_sochi/team/"
$arr[0,10] = "2025-12-22T04:41:10.631777+00:00"
$rng.Value = $arr

$rng = $ws.Range("A24:K24")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СОЧ"
$arr[0,1] = "ХК Сочи"
$arr[0,2] = "hc_sochi"
$arr[0,3] = "Швырёв Игорь"
$arr[0,4] = "98"
$arr[0,5] = "нападающий"
$arr[0,6] = "23300"
$arr[0,7] = "1369_СОЧ_швыревигорь"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/hc_sochi/team/"
$arr[0,10] = "2025-12-22T04:41:10.631800+00:00"
$rng.Value = $arr

$rng = $ws.Range("A25:K25")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СПР"
$arr[0,1] = "Спартак"
$arr[0,2] = "spartak"
$arr[0,3] = "Вишневский Дмитрий"
$arr[0,4] = "55"
$arr[0,5] = "защитник"
$arr[0,6] = "15299"
$arr[0,7] = "1369_СПР_вишневскийдмитрий"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/spartak/team/"
$arr[0,10] = "2025-12-22T04:41:12.277636+00:00"
$rng.Value = $arr

$rng = $ws.Range("A26:K26")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СПР"
$arr[0,1] = "Спартак"
$arr[0,2] = "spartak"
$arr[0,3] = "Коростелёв Никита"
$arr[0,4] = "35"
$arr[0,5] = "нападающий"
$arr[0,6] = "22149"
$arr[0,7] = "1369_СПР_коростелевникита"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/spartak/team/"
$arr[0,10] = "2025-12-22T04:41:12.277664+00:00"
$rng.Value = $arr

$rng = $ws.Range("A27:K27")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СПР"
$arr[0,1] = "Спартак"
$arr[0,2] = "spartak"
$arr[0,3] = "Порядин Павел"
$arr[0,4] = "24"
$arr[0,5] = "нападающий"
$arr[0,6] = "19258"
$arr[0,7] = "1369_СПР_порядинпавел"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/spartak/team/"
$arr[0,10] = "2025-12-22T04:41:12.277681+00:00"
$rng.Value = $arr

$rng = $ws.Range("A28:K28")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "СЮЛ"
$arr[0,1] = "Салават Юлаев"
$arr[0,2] = "salavat_yulaev"
$arr[0,3] = "Алалыкин Данил"
$arr[0,4] = "61"
$arr[0,5] = "нападающий"
$arr[0,6] = "34493"
$arr[0,7] = "1369_СЮЛ_алалыкинданил"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/salavat_yulaev/team/"
$arr[0,10] = "2025-12-22T04:41:13.930041+00:00"
$rng.Value = $arr

$rng = $ws.Range("A29:K29")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "ТОР"
$arr[0,1] = "Торпедо"
$arr[0,2] = "torpedo"
$arr[0,3] = "Кручинин Алексей"
$arr[0,4] = "78"
$arr[0,5] = "нападающий"
$arr[0,6] = "16355"
$arr[0,7] = "1369_ТОР_кручининалексей"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/torpedo/team/"
$arr[0,10] = "2025-12-22T04:41:16.113068+00:00"
$rng.Value = $arr

$rng = $ws.Range("A30:K30")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "ТОР"
$arr[0,1] = "Торпедо"
$arr[0,2] = "torpedo"
$arr[0,3] = "Принс Шэйн"
$arr[0,4] = "18"
$arr[0,5] = "нападающий"
$arr[0,6] = "19060"
$arr[0,7] = "1369_ТОР_принсшэйн"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/torpedo/team/"
$arr[0,10] = "2025-12-22T04:41:16.113092+00:00"
$rng.Value = $arr

$rng = $ws.Range("A31:K31")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = "ТОР"
$arr[0,1] = "Торпедо"
$arr[0,2] = "torpedo"
$arr[0,3] = "Шавин Никита"
$arr[0,4] = "80"
$arr[0,5] = "нападающий"
$arr[0,6] = "30548"
$arr[0,7] = "1369_ТОР_шавинникита"
$arr[0,8] = "injured_active"
$arr[0,9] = "https://www.khl.ru/clubs/torpedo/team/"
$arr[0,10] = "2025-12-22T04:41:16.113100+00:00"
$rng.Value = $arr

$c = $ws.Range("K32")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:41:19.988597+00:00"

$c = $ws.Range("K33")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:41:19.988623+00:00"

$c = $ws.Range("K34")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:41:22.056213+00:00"

$c = $ws.Range("K35")
$c.NumberFormat = "@"
$c.Value = "2025-12-22T04:41:22.056240+00:00"


$ws = $wb.Worksheets.Item("returned")

$rng = $ws.Range("A2:G2")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = "БАР"
$arr[0,1] = "Барыс"
$arr[0,2] = "Уотерспун Тайлер"
$arr[0,3] = "1369_БАР_уотерспунтайлер"
$arr[0,4] = "RETURN"
$arr[0,5] = "2025-12-22T12:41:22.562459+08:00"
$arr[0,6] = "2025-12-22"
$rng.Value = $arr

$rng = $ws.Range("A3:G3")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = "СИБ"
$arr[0,1] = "Сибирь"
$arr[0,2] = "Першаков Александр"
$arr[0,3] = "1369_СИБ_першаковалександр"
$arr[0,4] = "RETURN"
$arr[0,5] = "2025-12-22T12:41:22.562459+08:00"
$arr[0,6] = "2025-12-22"
$rng.Value = $arr


$ws = $wb.Worksheets.Item("new_injured")

$rng = $ws.Range("A2:G2")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = "СИБ"
$arr[0,1] = "Сибирь"
$arr[0,2] = "Люзенков Илья"
$arr[0,3] = "1369_СИБ_люзенковилья"
$arr[0,4] = "INJURED_NEW"
$arr[0,5] = "2025-12-22T12:41:22.562459+08:00"
$arr[0,6] = "2025-12-22"
$rng.Value = $arr

$rng = $ws.Range("A3:G3")
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = "ТОР"
$arr[0,1] = "Торпедо"
$arr[0,2] = "Шавин Никита"
$arr[0,3] = "1369_ТОР_шавинникита"
$arr[0,4] = "INJURED_NEW"
$arr[0,5] = "2025-12-22T12:41:22.562459+08:00"
$arr[0,6] = "2025-12-22"
$rng.Value = $arr


Write-Host "Applied KHL injuries update for 2025-12-22 run."